$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so Excel keeps them as text
$textCells = @("D4", "D5", "D6", "D8", "D10", "D13", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D30", "D32", "D35", "D36", "D40", "D45", "D47", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '68.059.50'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '3.791.89'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = '601.01'
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").Value = '165.46'
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.518'
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("E9").Value = '  -0.75%  '
$ws.Range("D10").Value = '0.452'
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("E11").Value = '  +2.58%  '
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("D13").Value = '35.78'
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("D14").Value = '4.424.09'
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("D15").Value = '3.778.98'
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("D16").Value = '68.037.27'
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("D17").Value = '18.37'
$ws.Range("E17").Value = '  -0.92%  '
$ws.Range("E18").Value = '  +1.91%  '
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").Value = '461.37'
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").Value = '9.73'
$ws.Range("E21").Value = '  -1.55%  '
$ws.Range("D22").Value = '0.698'
$ws.Range("D23").Value = '0.0000150'
$ws.Range("E23").Value = '  -2.32%  '
$ws.Range("D24").Value = '82.80'
$ws.Range("E24").Value = '  -0.86%  '
$ws.Range("D25").Value = '12.06'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("E26").Value = '  +0.38%  '
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").Value = '3.940.03'
$ws.Range("E29").Value = '  -0.21%  '
$ws.Range("D30").Value = '7.40'
$ws.Range("E30").Value = '  +2.27%  '
$ws.Range("E31").Value = '  -5.63%  '
$ws.Range("D32").Value = '2.22'
$ws.Range("E32").Value = '  -1.62%  '
$ws.Range("E33").Value = '  -1.05%  '
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D35").Value = '9.00'
$ws.Range("E35").Value = '  -0.80%  '
$ws.Range("D36").Value = '0.1000'
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("E38").Value = '  -3.29%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").Value = '0.988'
$ws.Range("E40").Value = '  -0.67%  '
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("E44").Value = '  -1.49%  '
$ws.Range("D45").Value = '43.00'
$ws.Range("E45").Value = '  -2.46%  '
$ws.Range("E46").Value = '  +0.25%  '
$ws.Range("D47").Value = '8.36'
$ws.Range("E47").Value = '  +0.48%  '
$ws.Range("E48").Value = '  +2.87%  '
$ws.Range("E49").Value = '  +7.06%  '
$ws.Range("D50").Value = '392.99'
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("D51").Value = '26.79'
$ws.Range("E51").Value = '  +1.80%  '

# Reset style index on text-forced cells back to default (keep text type, drop number format)
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
